# Finish the classification ("reclassificacao") for the remaining tweets in
# the "Teste" sheet: rows 72-201 already had their tweet text (column A)
# filled in but were missing the classification score in column B. This
# fills in those scores, completing the column B classification for every
# remaining row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teste")

$classifications = @{
    72 = 2
    73 = 3
    74 = 2
    75 = 3
    76 = 3
    77 = 3
    78 = 1
    79 = 2
    80 = 1
    81 = 1
    82 = 2
    83 = 3
    84 = 0
    85 = 0
    86 = 2
    87 = 2
    88 = 0
    89 = 0
    90 = 2
    91 = 0
    92 = 3
    93 = 0
    94 = 0
    95 = 2
    96 = 3
    97 = 2
    98 = 1
    99 = 1
    100 = 0
    101 = 3
    102 = 2
    103 = 2
    104 = 3
    105 = 3
    106 = 2
    107 = 2
    108 = 3
    109 = 1
    110 = 3
    111 = 2
    112 = 2
    113 = 1
    114 = 2
    115 = 3
    116 = 2
    117 = 3
    118 = 1
    119 = 3
    120 = 2
    121 = 2
    122 = 1
    123 = 0
    124 = 2
    125 = 2
    126 = 1
    127 = 0
    128 = 3
    129 = 2
    130 = 3
    131 = 3
    132 = 2
    133 = 3
    134 = 3
    135 = 3
    136 = 3
    137 = 1
    138 = 2
    139 = 1
    140 = 3
    141 = 3
    142 = 3
    143 = 0
    144 = 3
    145 = 3
    146 = 2
    147 = 0
    148 = 3
    149 = 1
    150 = 3
    151 = 3
    152 = 2
    153 = 2
    154 = 3
    155 = 2
    156 = 2
    157 = 2
    158 = 3
    159 = 1
    160 = 3
    161 = 3
    162 = 2
    163 = 3
    164 = 3
    165 = 2
    166 = 3
    167 = 3
    168 = 3
    169 = 2
    170 = 3
    171 = 3
    172 = 3
    173 = 2
    174 = 2
    175 = 3
    176 = 2
    177 = 3
    178 = 2
    179 = 0
    180 = 3
    181 = 3
    182 = 1
    183 = 2
    184 = 3
    185 = 3
    186 = 3
    187 = 2
    188 = 3
    189 = 2
    190 = 2
    191 = 0
    192 = 2
    193 = 2
    194 = 3
    195 = 3
    196 = 2
    197 = 2
    198 = 3
    199 = 2
    200 = 2
    201 = 1
}

foreach ($row in ($classifications.Keys | Sort-Object {[int]$_})) {
    $ws.Cells.Item([int]$row, 2).Value = $classifications[$row]
}

# Reflect the editing cursor ending on the last newly classified cell.
$ws.Range("B202").Select()
